$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("图片路径")

$ws.Range("A2").Value = "/usr/local/picture/1.jpg"
$ws.Range("A3").Value = "/usr/local/picture/2.jpeg"
$ws.Range("A4").Value = "/usr/local/picture/3.png"
$ws.Range("A5").Value = "/usr/local/picture/4.gif"
$ws.Range("A6").Value = "/usr/local/picture/5.jfif"

$ws.Range("B5").Select()
